# ---------------------------------------------------------------------------
# Add "2022-Q3" quarterly snapshot to the 601766-中国中车 holders workbook.
#
#  1. Insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q2"),
#     populated with the new quarter's fund-holder table.
#  2. Insert a new row into "总计" for the "2022-Q3" summary line, pushing the
#     older quarters down by one row.
#
# Helper: Set-TextCell forces a value to be stored as TEXT (so fund codes
# like "004497"/"015678" keep their leading zeros, and numeric-looking
# figures like "3.40" keep their trailing zero) while keeping the cell's
# style index unchanged (no stray NumberFormat/quote-prefix styling left
# behind) by resetting to the built-in "Normal" style right after the
# write.
# ---------------------------------------------------------------------------

function Set-TextCell($cell, $val) {
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2022-Q2")

$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Bring over the header-row formatting (bold + border, B1:H1) and the
# column-A row-index formatting (bold + border) from the existing
# "2022-Q2" sheet, which uses the same layout.
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q2.Range("A2").Copy()
$q3.Range("A2:A13").PasteSpecial(-4122)

# Header row
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# Fund holder rows (row 2 .. row 13), columns:
#  A=index(0-based,numeric)  B=code(text)  C=name(text)  D=scale(text)
#  E=position%(text)  F=weight%(text)  G=marketvalue(text)  H=rank(numeric)
$rows = @(
    @(0,  "515900", "博时中证央企创新驱动ETF",       "36.57", "98.62", "3.05",  "1.1154", 6),
    @(1,  "516950", "银华中证基建ETF",               "11.07", "97.93", "7.72",  "0.8546", 2),
    @(2,  "515600", "广发中证央企创新驱动ETF",       "14.84", "98.78", "3.04",  "0.4511", 6),
    @(3,  "515680", "嘉实中证央企创新驱动ETF",       "14.64", "99.23", "3.04",  "0.4451", 6),
    @(4,  "159635", "华夏中证基建ETF",               "3.40",  "99.03", "7.81",  "0.2655", 2),
    @(5,  "160135", "南方中证高铁产业指数（LOF）",   "1.84",  "95.01", "13.96", "0.2569", 2),
    @(6,  "159619", "国泰中证基建ETF",               "3.30",  "98.76", "7.66",  "0.2528", 2),
    @(7,  "159974", "富国中证央企创新驱动ETF",       "4.89",  "99.47", "3.05",  "0.1491", 6),
    @(8,  "160639", "鹏华中证高铁产业指数（LOF）A",  "0.75",  "94.62", "13.85", "0.1039", 2),
    @(9,  "004497", "前海开源多元策略灵活配置混合C", "1.68",  "93.04", "4.88",  "0.0820", 4),
    @(10, "004496", "前海开源多元策略灵活配置混合A", "0.91",  "93.04", "4.88",  "0.0444", 4),
    @(11, "015678", "鹏华中证高铁产业指数（LOF）C",  "0.06",  "94.62", "13.85", "0.0083", 2)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r,1).Value = $row[0]
    Set-TextCell $q3.Cells.Item($r,2) $row[1]
    $q3.Cells.Item($r,3).Value = $row[2]
    Set-TextCell $q3.Cells.Item($r,4) $row[3]
    Set-TextCell $q3.Cells.Item($r,5) $row[4]
    Set-TextCell $q3.Cells.Item($r,6) $row[5]
    Set-TextCell $q3.Cells.Item($r,7) $row[6]
    $q3.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Step 2: insert the "2022-Q3" summary row into "总计"
# ---------------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()
$zongji.Range("B2:D2").Style = "Normal"
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)

$zongji.Cells.Item(2,1).Value = 0
$zongji.Cells.Item(2,2).Value = "2022-Q3"
$zongji.Cells.Item(2,3).Value = 12
$zongji.Cells.Item(2,4).Value = 4.03
